$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.071.25'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '3.172.05'
$ws.Range("E3").Value = '  -4.57%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'591.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").Value = "'134.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.45%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.167.65'
$ws.Range("E8").Value = '  -4.62%  '
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -5.52%  '
$ws.Range("D11").Value = "'5.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.99%  '
$ws.Range("E12").Value = '  -3.20%  '
$ws.Range("E13").Value = '  -4.06%  '
$ws.Range("D14").Value = "'34.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '3.691.83'
$ws.Range("E15").Value = '  -4.59%  '
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '3.172.22'
$ws.Range("E17").Value = '  -4.48%  '
$ws.Range("D18").Value = '63.039.06'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  -4.03%  '
$ws.Range("D20").Value = "'461.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.05%  '
$ws.Range("D21").Value = "'13.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("E22").Value = '  -5.16%  '
$ws.Range("D23").Value = "'7.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("D25").Value = "'83.28"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E28").Value = '  -3.74%  '
$ws.Range("D29").Value = "'7.74"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.65%  '
$ws.Range("D30").Value = "'6.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.82%  '
$ws.Range("E31").Value = '  -5.96%  '
$ws.Range("D32").Value = "'27.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.17%  '
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("E34").Value = '  -6.49%  '
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("E36").Value = '  -4.04%  '
$ws.Range("D37").Value = "'51.47"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").Value = '0.0₃0709'
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("D40").Value = "'405.27"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.80%  '
$ws.Range("E41").Value = '  -2.52%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("D43").Value = '2.816.07'
$ws.Range("E43").Value = '  -9.07%  '
$ws.Range("D44").Value = "'2.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.75%  '
$ws.Range("E45").Value = '  -6.04%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = "'2.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.28%  '
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Value = "'123.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").Value = "'34.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.04%  '
